$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive original leading text.
$paras = $d.Paragraphs
$count = $paras.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Help developing features, deploying")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "ERROR: target paragraph not found"
} else {
    $full = $target.Range
    # Exclude the trailing paragraph mark from the replaced range so the
    # paragraph itself (and its mark) are preserved, only its runs change.
    $contentRange = $d.Range($full.Start, $full.End - 1)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="32DE8A7C" w14:textId="621AA40E" w:rsidR="0064141A" w:rsidRPr="0064141A" w:rsidRDefault="0064141A" w:rsidP="0064141A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:jc w:val="both"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>D</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>evelop</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>ed</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>MLOps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve">features, deploying </w:t></w:r><w:r w:rsidRPr="006A470A"><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>and maintaining</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> ML </w:t></w:r><w:r w:rsidRPr="006A470A"><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>services</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> Document Processing and Extraction services</w:t></w:r><w:r w:rsidRPr="006A470A"><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> using Python and </w:t></w:r><w:r w:rsidRPr="006A470A"><w:rPr><w:b/><w:bCs/><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>Rust</w:t></w:r><w:r w:rsidRPr="006A470A"><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidRPr="0064141A"><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    [void]$contentRange.InsertXML($xml)
    Write-Host "Paragraph runs replaced."
}
